# RevC_PowerCalculations.xlsx -- "Bulk of schematic work finished, beginning layout"
#
# Adds a Min-Efficiency requirement row, a Results: block with Pout/Pin/
# max-input-current calculations, relocates the "LT8631" note further down
# the sheet, widens column A to fit the longer labels, and nudges the
# window/selection state to reflect where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window state: minimize the workbook window -----------------------
$wb.Windows.Item(1).WindowState = -4140   # xlMinimized

# --- Column A: widen to fit the new, longer row labels -----------------
$ws.Columns.Item(1).ColumnWidth = 27.5

# --- New requirement row: Min Efficiency --------------------------------
$ws.Range("A8").Value = "Min Efficiency "
$ws.Range("B8").Value = ">"
$ws.Range("C8").Value = 80
$ws.Range("D8").Value = "%"

# --- Relocate the "LT8631" note from row 9 down to row 21 --------------
$ws.Range("A9").Cut($ws.Range("A21"))

# --- Results: block -------------------------------------------------
$ws.Range("A12").Value = "Results:"
$ws.Range("A12").Font.Bold = $true

$ws.Range("A13").Value = "Pout"
$ws.Range("B13").Formula = "=C7/1000*C4"

$ws.Range("A14").Value = "Pin"
$ws.Range("B14").Formula = "=B13*(1/C8)*100"

$ws.Range("A15").Value = "Current In Max @ 18v"
$ws.Range("B15").Formula = "=B14/18"

$ws.Range("A16").Value = "Current In Max @ 50v"
$ws.Range("B16").Formula = "=B14/50"

# --- Leave the selection where the author left off ----------------------
$ws.Range("E13").Select()
